$wb = $excel.ActiveWorkbook

# Sheet: 展览 (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1790
$ws1.Range("F4").Value = 454
$ws1.Range("F9").Value = 1738
$ws1.Range("F10").Value = 368
$ws1.Range("F12").Value = 811
$ws1.Range("F13").Value = 339
$ws1.Range("F15").Value = 12822
$ws1.Range("F16").Value = 12813
$ws1.Range("F17").Value = 959
$ws1.Range("F22").Value = 569
$ws1.Range("F23").Value = 2010
$ws1.Range("F24").Value = 31
$ws1.Range("F27").Value = 43

# Sheet: 演出 (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value = 10

# Sheet: 本地生活 (sheet3)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 167

# Sheet: 全部类型 (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 167
$ws4.Range("F5").Value = 1790
$ws4.Range("F6").Value = 454
$ws4.Range("F14").Value = 1738
$ws4.Range("F15").Value = 368
$ws4.Range("F17").Value = 811
$ws4.Range("F18").Value = 339
$ws4.Range("F21").Value = 12822
$ws4.Range("F22").Value = 12813
$ws4.Range("F23").Value = 959
$ws4.Range("F28").Value = 569
$ws4.Range("F30").Value = 10
$ws4.Range("F31").Value = 2010
$ws4.Range("F32").Value = 31
$ws4.Range("F37").Value = 43
